# Apply updated cryptocurrency data (price & 1h volume change) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '22.022.61'
Set-TextCell $ws.Range("E2") '  -1.91%  '

# Row 3
Set-TextCell $ws.Range("D3") '1.554.58'
Set-TextCell $ws.Range("E3") '  -1.18%  '

# Row 5
Set-TextCell $ws.Range("D5") '1.001'
Set-TextCell $ws.Range("E5") '  +0.00%  '

# Row 6
Set-TextCell $ws.Range("D6") '286.55'
Set-TextCell $ws.Range("E6") '  -0.48%  '

# Row 7
Set-TextCell $ws.Range("D7") '0.3763'
Set-TextCell $ws.Range("E7") '  +1.64%  '

# Row 8
Set-TextCell $ws.Range("D8") '0.3238'
Set-TextCell $ws.Range("E8") '  -2.48%  '

# Row 9
Set-TextCell $ws.Range("D9") '1.124'
Set-TextCell $ws.Range("E9") '  -2.49%  '

# Row 10
Set-TextCell $ws.Range("D10") '41.18'
Set-TextCell $ws.Range("E10") '  -12.78%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.07306'
Set-TextCell $ws.Range("E11") '  -2.58%  '

# Row 12
Set-TextCell $ws.Range("E12") '  -0.01%  '

# Row 13
Set-TextCell $ws.Range("D13") '19.51'
Set-TextCell $ws.Range("E13") '  -6.11%  '

# Row 14
Set-TextCell $ws.Range("D14") '5.713'
Set-TextCell $ws.Range("E14") '  -3.76%  '

# Row 15
Set-TextCell $ws.Range("D15") '6.831'
Set-TextCell $ws.Range("E15") '  -1.34%  '

# Row 16
Set-TextCell $ws.Range("D16") '1.551.31'
Set-TextCell $ws.Range("E16") '  -0.66%  '

# Row 17
Set-TextCell $ws.Range("E17") '  -3.11%  '

# Row 18
Set-TextCell $ws.Range("D18") '0.06647'
Set-TextCell $ws.Range("E18") '  -1.10%  '

# Row 19
Set-TextCell $ws.Range("D19") '85.09'
Set-TextCell $ws.Range("E19") '  -3.77%  '

# Row 20
Set-TextCell $ws.Range("D20") '6.431'
Set-TextCell $ws.Range("E20") '  +0.65%  '

# Row 21
Set-TextCell $ws.Range("D21") '1.000'
Set-TextCell $ws.Range("E21") '  +0.05%  '

# Row 22
Set-TextCell $ws.Range("D22") '16.00'
Set-TextCell $ws.Range("E22") '  -2.97%  '

# Row 23
Set-TextCell $ws.Range("D23") '11.56'
Set-TextCell $ws.Range("E23") '  -3.60%  '

# Row 24
Set-TextCell $ws.Range("D24") '22.032.27'
Set-TextCell $ws.Range("E24") '  -1.79%  '

# Row 25
Set-TextCell $ws.Range("D25") '2.245'
Set-TextCell $ws.Range("E25") '  -5.83%  '

# Row 26
Set-TextCell $ws.Range("D26") '2.532'
Set-TextCell $ws.Range("E26") '  -4.06%  '

# Row 27
Set-TextCell $ws.Range("D27") '150.26'
Set-TextCell $ws.Range("E27") '  -0.32%  '

# Row 28
Set-TextCell $ws.Range("D28") '18.89'
Set-TextCell $ws.Range("E28") '  -3.62%  '

# Row 29
Set-TextCell $ws.Range("D29") '4.834'
Set-TextCell $ws.Range("E29") '  -2.66%  '

# Row 30
Set-TextCell $ws.Range("D30") '1.724.59'
Set-TextCell $ws.Range("E30") '  -1.00%  '

# Row 31
Set-TextCell $ws.Range("D31") '120.22'
Set-TextCell $ws.Range("E31") '  -3.83%  '

# Row 32
Set-TextCell $ws.Range("D32") '1.118'
Set-TextCell $ws.Range("E32") '  +2.11%  '

# Row 33
Set-TextCell $ws.Range("D33") '5.929'
Set-TextCell $ws.Range("E33") '  -2.54%  '

# Row 34
Set-TextCell $ws.Range("D34") '1.662'
Set-TextCell $ws.Range("E34") '  -16.30%  '

# Row 35
Set-TextCell $ws.Range("D35") '9.279'
Set-TextCell $ws.Range("E35") '  -6.21%  '

# Row 36
Set-TextCell $ws.Range("D36") '0.08161'
Set-TextCell $ws.Range("E36") '  -2.22%  '

# Row 37
Set-TextCell $ws.Range("B37") 'InternetComputer(DFINITY)'
Set-TextCell $ws.Range("C37") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range("D37") '5.230'
Set-TextCell $ws.Range("E37") '  -2.02%  '

# Row 38
Set-TextCell $ws.Range("B38") 'VeChain'
Set-TextCell $ws.Range("C38") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D38") '0.02287'
Set-TextCell $ws.Range("E38") '  -6.52%  '

# Row 39
Set-TextCell $ws.Range("D39") '0.06156'
Set-TextCell $ws.Range("E39") '  -3.69%  '

# Row 40
Set-TextCell $ws.Range("D40") '0.2117'
Set-TextCell $ws.Range("E40") '  -4.58%  '

# Row 41
Set-TextCell $ws.Range("E41") '  -7.02%  '

# Row 42
Set-TextCell $ws.Range("E42") '  -4.30%  '

# Row 43
Set-TextCell $ws.Range("D43") '1.000'
Set-TextCell $ws.Range("E43") '  +0.02%  '

# Row 44
Set-TextCell $ws.Range("D44") '0.5944'
Set-TextCell $ws.Range("E44") '  -4.79%  '

# Row 45
Set-TextCell $ws.Range("D45") '13.58'
Set-TextCell $ws.Range("E45") '  -3.21%  '

# Row 46
Set-TextCell $ws.Range("D46") '3.725'
Set-TextCell $ws.Range("E46") '  -1.33%  '

# Row 47
Set-TextCell $ws.Range("E47") '  -5.40%  '

# Row 48
Set-TextCell $ws.Range("D48") '1.948'
Set-TextCell $ws.Range("E48") '  -4.81%  '

# Row 49
Set-TextCell $ws.Range("D49") '120.06'
Set-TextCell $ws.Range("E49") '  -3.83%  '

# Row 50
Set-TextCell $ws.Range("D50") '1.155'
Set-TextCell $ws.Range("E50") '  -4.35%  '

# Row 51
Set-TextCell $ws.Range("D51") '0.06925'
Set-TextCell $ws.Range("E51") '  -3.77%  '
